# Bugfixed the naive forecaster component module
#
# The underlying forecast-generation run shifted by one period (the very
# first observation date/year pair was dropped and every later row moved
# up by one), and the AR(2) forecast coefficients (columns C and E) were
# recomputed with the corrected data. This script reproduces both effects
# on the existing worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the obsolete first data row (old row 2); this shifts every
#    remaining row up by one and naturally realigns columns A (date),
#    B (y_0 year) and D (y_1 year), which line up exactly with the old
#    row below them.
$ws.Rows.Item(2).Delete()

# 2) Recompute / refresh the forecast columns C (y_0_forecast) and
#    E (y_1_forecast) with the corrected values. Some rows have no
#    forecast (blank).
$C_VALUES = @($null,$null,$null,-1.317619545389281,$null,2.771597318554297,$null,1.799362536952542,$null,2.123182427147152,4.880442637054072,7.317297369612819,5.941867202078877,4.260319658857736,0.292749233164491,-0.05262415810141086,2.032207428223742,1.459778471779982,2.775332754349846,3.002208343813528,3.2651197821016,3.565025829754953,4.01493878081518,4.020433260014977,3.283136334808323,3.444206290325491,3.479628752085517,3.53224976671227,2.31260691849986,2.667234932970275,-0.985458715495402,-0.985458715495402,-4.853362183897836,-4.511102905979703,-4.365687260408224,-4.365687260408224,2.700663803921799,1.386772772629241,1.27347919322387,1.27347919322387,-0.185315122156382,-0.9537175292835154,-1.339436245206127,-1.339436245206127,-2.5174493871855,-3.303819519576723,-3.451527003230626,-3.451527003230626,-1.652703512303566,-1.704805397136089,-1.75044229618867)

$E_VALUES = @($null,$null,$null,$null,$null,$null,$null,$null,$null,4.792854588620821,1.641301872652501,7.239454936865775,2.672847571394987,0.4361429468412448,2.181874035977249,1.392195163617171,2.13692496326825,1.493220091771108,2.200426660963761,2.622364272988187,2.271936475508851,2.446228176258058,3.091110147865495,3.113086948791377,2.403408536719187,2.496958452261078,2.54748094003614,2.757421718286168,2.363182008239928,2.567662999186382,-1.999977666418695,-1.373617952268746,0.6599234717970859,0.9582724917052587,1.582150300418306,1.324283050325015,2.813733698830267,0.8813242377093244,0.4287685802702779,-0.5025420863900898,1.715279842342743,0.9049225073274991,-0.6598239038080322,-2.371854438773213,1.089612584330668,-0.3655818470008065,-1.070961900287937,-1.436963918858969,0.4469810487905734,0.1460701281005727,-0.1564272439962933)

$startRow = 2
for ($i = 0; $i -lt $C_VALUES.Length; $i++) {
    $row = $startRow + $i

    $cCell = $ws.Cells.Item($row, 3)
    if ($null -eq $C_VALUES[$i]) {
        $cCell.ClearContents()
    } else {
        $cCell.Value = $C_VALUES[$i]
    }

    $eCell = $ws.Cells.Item($row, 5)
    if ($null -eq $E_VALUES[$i]) {
        $eCell.ClearContents()
    } else {
        $eCell.Value = $E_VALUES[$i]
    }
}

Write-Output "Applied naive forecaster bugfix to $($C_VALUES.Length) rows"
